$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D): force Text format so numeric-looking strings
# (e.g. '212.75', '143.30') are preserved exactly, matching the
# original inline-string cell content instead of being coerced
# into floating point numbers by Excel's automatic type detection.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.290.60"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.75"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0614"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.45"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.832.63"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.603.66"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.02"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.265.93"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.15"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "201.65"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.30"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.57"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.161.48"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.787"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.37"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.743.87"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.90"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.05"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0507"

# Volume(1h) column (E): plain percentage strings, never
# numeric-looking so no coercion risk.
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +0.66%  "
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("E10").Value = "  +2.17%  "
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("E14").Value = "  +0.39%  "
$ws.Range("E16").Value = "  +0.46%  "
$ws.Range("E17").Value = "  +2.74%  "
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("E20").Value = "  -1.18%  "
$ws.Range("E21").Value = "  +0.68%  "
$ws.Range("E22").Value = "  +0.44%  "
$ws.Range("E23").Value = "  +0.42%  "
$ws.Range("E24").Value = "  +0.28%  "
$ws.Range("E25").Value = "  +1.21%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  -1.50%  "
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("E29").Value = "  +2.19%  "
$ws.Range("E30").Value = "  +5.47%  "
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("E32").Value = "  +2.73%  "
$ws.Range("E33").Value = "  -1.36%  "
$ws.Range("E34").Value = "  +1.12%  "
$ws.Range("E35").Value = "  +1.44%  "
$ws.Range("E36").Value = "  +3.38%  "
$ws.Range("E37").Value = "  +1.39%  "
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("E39").Value = "  +0.92%  "
$ws.Range("E40").Value = "  +0.45%  "
$ws.Range("E41").Value = "  +0.87%  "
$ws.Range("E42").Value = "  +4.08%  "
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("E44").Value = "  +0.21%  "
$ws.Range("E45").Value = "  -0.99%  "
$ws.Range("E46").Value = "  +14.40%  "
$ws.Range("E47").Value = "  +0.92%  "
$ws.Range("E48").Value = "  +0.90%  "
$ws.Range("E49").Value = "  +0.43%  "
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("E51").Value = "  +0.04%  "
